$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; this shifts rows 45-171 down to 46-172.
$ws.Rows("45:45").Insert()

# Populate the newly inserted row 45 with the new data record.
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = 45044
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100108
$ws.Range("H45").Value = "Tropicales y subtropicales"
$ws.Range("I45").Value = 100108003
$ws.Range("J45").Value = "Maracuyá"
$ws.Range("K45").Value = "Sin especificar"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 70
$ws.Range("N45").Value = 34000
$ws.Range("O45").Value = 35000
$ws.Range("P45").Value = 34429
$ws.Range("Q45").Value = "$/caja 20 kilos"
$ws.Range("R45").Value = "Región de Arica y Parinacota"
$ws.Range("S45").Value = 1721
$ws.Range("T45").Value = 20
